$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price record was added to the dataset. In the source table this
# corresponds to inserting a new row at position 436, which pushes all the
# existing records (old rows 436-506) down by one (to rows 437-507).
$ws.Rows(436).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(436, 1).Value = 10
$ws.Cells.Item(436, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(436, 3).Value = "La Araucanía"
$ws.Cells.Item(436, 4).Value = 44951
$ws.Cells.Item(436, 5).Value = 9
$ws.Cells.Item(436, 6).Value = 100112040
$ws.Cells.Item(436, 7).Value = "Cilantro"
$ws.Cells.Item(436, 8).Value = "Sin especificar"
$ws.Cells.Item(436, 9).Value = "Primera"
$ws.Cells.Item(436, 10).Value = 55
$ws.Cells.Item(436, 11).Value = 6000
$ws.Cells.Item(436, 12).Value = 6000
$ws.Cells.Item(436, 13).Value = 6000
$ws.Cells.Item(436, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(436, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(436, 16).Value = 3000
$ws.Cells.Item(436, 17).Value = 2
$ws.Cells.Item(436, 18).Value = "Hortaliza"
